$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("float transitions")

# --- Build a stable style palette in a scratch area (column K) sourced from
# --- original cells, BEFORE any of the table values/styles are touched.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("K1").PasteSpecial(-4122) | Out-Null
$ws.Range("B2").Copy() | Out-Null
$ws.Range("K2").PasteSpecial(-4122) | Out-Null
$ws.Range("C2").Copy() | Out-Null
$ws.Range("K3").PasteSpecial(-4122) | Out-Null
$ws.Range("D2").Copy() | Out-Null
$ws.Range("K4").PasteSpecial(-4122) | Out-Null
$ws.Range("B3").Copy() | Out-Null
$ws.Range("K5").PasteSpecial(-4122) | Out-Null
$ws.Range("I3").Copy() | Out-Null
$ws.Range("K6").PasteSpecial(-4122) | Out-Null

# --- Seed new shared strings in the exact order required: S_FL_NUMQ, S_FL_DECQ, S_FL_DEC1
$ws.Range("H3").Value = "S_FL_NUMQ"
$ws.Range("G2").Value = "S_FL_DECQ"
$ws.Range("G3").Value = "S_FL_DEC1"

# --- Apply final styles + values to every cell in the transition table (rows 2-12).
# row 2
$ws.Range("K1").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null
$ws.Range("A2").Value = 0
$ws.Range("K2").Copy() | Out-Null
$ws.Range("B2").PasteSpecial(-4122) | Out-Null
$ws.Range("B2").Value = "S_FL_START"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("C2").PasteSpecial(-4122) | Out-Null
$ws.Range("C2").Value = "T_FL_ERROR"
$ws.Range("K4").Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4122) | Out-Null
$ws.Range("D2").Value = "S_FL_START"
$ws.Range("K4").Copy() | Out-Null
$ws.Range("E2").PasteSpecial(-4122) | Out-Null
$ws.Range("E2").Value = "S_FL_NUM"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("F2").PasteSpecial(-4122) | Out-Null
$ws.Range("F2").Value = "T_FL_ERROR"
$ws.Range("K4").Copy() | Out-Null
$ws.Range("G2").PasteSpecial(-4122) | Out-Null
$ws.Range("G2").Value = "S_FL_DECQ"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("H2").PasteSpecial(-4122) | Out-Null
$ws.Range("H2").Value = "T_FL_ERROR"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("I2").PasteSpecial(-4122) | Out-Null
$ws.Range("I2").Value = "T_FL_ERROR"

# row 3
$ws.Range("K1").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4122) | Out-Null
$ws.Range("A3").Value = 1
$ws.Range("K5").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4122) | Out-Null
$ws.Range("B3").Value = "S_FL_NUM"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("C3").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Value = "T_FL_ERROR"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("D3").PasteSpecial(-4122) | Out-Null
$ws.Range("D3").Value = "T_FL_ERROR"
$ws.Range("K4").Copy() | Out-Null
$ws.Range("E3").PasteSpecial(-4122) | Out-Null
$ws.Range("E3").Value = "S_FL_NUM"
$ws.Range("K4").Copy() | Out-Null
$ws.Range("F3").PasteSpecial(-4122) | Out-Null
$ws.Range("F3").Value = "S_FL_EXP"
$ws.Range("K4").Copy() | Out-Null
$ws.Range("G3").PasteSpecial(-4122) | Out-Null
$ws.Range("G3").Value = "S_FL_DEC1"
$ws.Range("K4").Copy() | Out-Null
$ws.Range("H3").PasteSpecial(-4122) | Out-Null
$ws.Range("H3").Value = "S_FL_NUMQ"
$ws.Range("K6").Copy() | Out-Null
$ws.Range("I3").PasteSpecial(-4122) | Out-Null
$ws.Range("I3").Value = "T_FL_FLOAT"

# row 4
$ws.Range("K1").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null
$ws.Range("A4").Value = 2
$ws.Range("K5").Copy() | Out-Null
$ws.Range("B4").PasteSpecial(-4122) | Out-Null
$ws.Range("B4").Value = "S_FL_NUMQ"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("C4").PasteSpecial(-4122) | Out-Null
$ws.Range("C4").Value = "T_FL_ERROR"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("D4").PasteSpecial(-4122) | Out-Null
$ws.Range("D4").Value = "T_FL_ERROR"
$ws.Range("K4").Copy() | Out-Null
$ws.Range("E4").PasteSpecial(-4122) | Out-Null
$ws.Range("E4").Value = "S_FL_NUM"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("F4").PasteSpecial(-4122) | Out-Null
$ws.Range("F4").Value = "T_FL_ERROR"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("G4").PasteSpecial(-4122) | Out-Null
$ws.Range("G4").Value = "T_FL_ERROR"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("H4").PasteSpecial(-4122) | Out-Null
$ws.Range("H4").Value = "T_FL_ERROR"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("I4").PasteSpecial(-4122) | Out-Null
$ws.Range("I4").Value = "T_FL_ERROR"

# row 5
$ws.Range("K1").Copy() | Out-Null
$ws.Range("A5").PasteSpecial(-4122) | Out-Null
$ws.Range("A5").Value = 3
$ws.Range("K5").Copy() | Out-Null
$ws.Range("B5").PasteSpecial(-4122) | Out-Null
$ws.Range("B5").Value = "S_FL_DEC"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("C5").PasteSpecial(-4122) | Out-Null
$ws.Range("C5").Value = "T_FL_ERROR"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4122) | Out-Null
$ws.Range("D5").Value = "T_FL_ERROR"
$ws.Range("K4").Copy() | Out-Null
$ws.Range("E5").PasteSpecial(-4122) | Out-Null
$ws.Range("E5").Value = "S_FL_DEC"
$ws.Range("K4").Copy() | Out-Null
$ws.Range("F5").PasteSpecial(-4122) | Out-Null
$ws.Range("F5").Value = "S_FL_EXP"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("G5").PasteSpecial(-4122) | Out-Null
$ws.Range("G5").Value = "T_FL_ERROR"
$ws.Range("K4").Copy() | Out-Null
$ws.Range("H5").PasteSpecial(-4122) | Out-Null
$ws.Range("H5").Value = "S_FL_DECQ"
$ws.Range("K6").Copy() | Out-Null
$ws.Range("I5").PasteSpecial(-4122) | Out-Null
$ws.Range("I5").Value = "T_FL_FLOAT"

# row 6
$ws.Range("K1").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4122) | Out-Null
$ws.Range("A6").Value = 4
$ws.Range("K5").Copy() | Out-Null
$ws.Range("B6").PasteSpecial(-4122) | Out-Null
$ws.Range("B6").Value = "S_FL_DECQ"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("C6").PasteSpecial(-4122) | Out-Null
$ws.Range("C6").Value = "T_FL_ERROR"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4122) | Out-Null
$ws.Range("D6").Value = "T_FL_ERROR"
$ws.Range("K4").Copy() | Out-Null
$ws.Range("E6").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").Value = "S_FL_DEC"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("F6").PasteSpecial(-4122) | Out-Null
$ws.Range("F6").Value = "T_FL_ERROR"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("G6").PasteSpecial(-4122) | Out-Null
$ws.Range("G6").Value = "T_FL_ERROR"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("H6").PasteSpecial(-4122) | Out-Null
$ws.Range("H6").Value = "T_FL_ERROR"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("I6").PasteSpecial(-4122) | Out-Null
$ws.Range("I6").Value = "T_FL_ERROR"

# row 7
$ws.Range("K1").Copy() | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null
$ws.Range("A7").Value = 5
$ws.Range("K5").Copy() | Out-Null
$ws.Range("B7").PasteSpecial(-4122) | Out-Null
$ws.Range("B7").Value = "S_FL_DEC1"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("C7").PasteSpecial(-4122) | Out-Null
$ws.Range("C7").Value = "T_FL_ERROR"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4122) | Out-Null
$ws.Range("D7").Value = "T_FL_ERROR"
$ws.Range("K4").Copy() | Out-Null
$ws.Range("E7").PasteSpecial(-4122) | Out-Null
$ws.Range("E7").Value = "S_FL_DEC"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("F7").PasteSpecial(-4122) | Out-Null
$ws.Range("F7").Value = "T_FL_ERROR"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("G7").PasteSpecial(-4122) | Out-Null
$ws.Range("G7").Value = "T_FL_ERROR"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("H7").PasteSpecial(-4122) | Out-Null
$ws.Range("H7").Value = "T_FL_ERROR"
$ws.Range("K6").Copy() | Out-Null
$ws.Range("I7").PasteSpecial(-4122) | Out-Null
$ws.Range("I7").Value = "T_FL_FLOAT"

# row 8
$ws.Range("K1").Copy() | Out-Null
$ws.Range("A8").PasteSpecial(-4122) | Out-Null
$ws.Range("A8").Value = 6
$ws.Range("K5").Copy() | Out-Null
$ws.Range("B8").PasteSpecial(-4122) | Out-Null
$ws.Range("B8").Value = "S_FL_EXP"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("C8").PasteSpecial(-4122) | Out-Null
$ws.Range("C8").Value = "T_FL_ERROR"
$ws.Range("K4").Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4122) | Out-Null
$ws.Range("D8").Value = "S_FL_EXPS"
$ws.Range("K4").Copy() | Out-Null
$ws.Range("E8").PasteSpecial(-4122) | Out-Null
$ws.Range("E8").Value = "S_FL_EXPD"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("F8").PasteSpecial(-4122) | Out-Null
$ws.Range("F8").Value = "T_FL_ERROR"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("G8").PasteSpecial(-4122) | Out-Null
$ws.Range("G8").Value = "T_FL_ERROR"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("H8").PasteSpecial(-4122) | Out-Null
$ws.Range("H8").Value = "T_FL_ERROR"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("I8").PasteSpecial(-4122) | Out-Null
$ws.Range("I8").Value = "T_FL_ERROR"

# row 9
$ws.Range("K1").Copy() | Out-Null
$ws.Range("A9").PasteSpecial(-4122) | Out-Null
$ws.Range("A9").Value = 7
$ws.Range("K5").Copy() | Out-Null
$ws.Range("B9").PasteSpecial(-4122) | Out-Null
$ws.Range("B9").Value = "S_FL_EXPS"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("C9").PasteSpecial(-4122) | Out-Null
$ws.Range("C9").Value = "T_FL_ERROR"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("D9").PasteSpecial(-4122) | Out-Null
$ws.Range("D9").Value = "T_FL_ERROR"
$ws.Range("K4").Copy() | Out-Null
$ws.Range("E9").PasteSpecial(-4122) | Out-Null
$ws.Range("E9").Value = "S_FL_EXPD"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("F9").PasteSpecial(-4122) | Out-Null
$ws.Range("F9").Value = "T_FL_ERROR"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("G9").PasteSpecial(-4122) | Out-Null
$ws.Range("G9").Value = "T_FL_ERROR"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("H9").PasteSpecial(-4122) | Out-Null
$ws.Range("H9").Value = "T_FL_ERROR"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("I9").PasteSpecial(-4122) | Out-Null
$ws.Range("I9").Value = "T_FL_ERROR"

# row 10
$ws.Range("K1").Copy() | Out-Null
$ws.Range("A10").PasteSpecial(-4122) | Out-Null
$ws.Range("A10").Value = 8
$ws.Range("K5").Copy() | Out-Null
$ws.Range("B10").PasteSpecial(-4122) | Out-Null
$ws.Range("B10").Value = "S_FL_EXPD"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("C10").PasteSpecial(-4122) | Out-Null
$ws.Range("C10").Value = "T_FL_ERROR"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4122) | Out-Null
$ws.Range("D10").Value = "T_FL_ERROR"
$ws.Range("K4").Copy() | Out-Null
$ws.Range("E10").PasteSpecial(-4122) | Out-Null
$ws.Range("E10").Value = "S_FL_EXPD"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("F10").PasteSpecial(-4122) | Out-Null
$ws.Range("F10").Value = "T_FL_ERROR"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("G10").PasteSpecial(-4122) | Out-Null
$ws.Range("G10").Value = "T_FL_ERROR"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("H10").PasteSpecial(-4122) | Out-Null
$ws.Range("H10").Value = "T_FL_ERROR"
$ws.Range("K6").Copy() | Out-Null
$ws.Range("I10").PasteSpecial(-4122) | Out-Null
$ws.Range("I10").Value = "T_FL_FLOAT"

# row 11
$ws.Range("K1").Copy() | Out-Null
$ws.Range("A11").PasteSpecial(-4122) | Out-Null
$ws.Range("A11").Value = 9
$ws.Range("K5").Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4122) | Out-Null
$ws.Range("B11").Value = "T_FL_FLOAT"
$ws.Range("K6").Copy() | Out-Null
$ws.Range("C11").PasteSpecial(-4122) | Out-Null
$ws.Range("C11").Value = "T_FL_FLOAT"
$ws.Range("K6").Copy() | Out-Null
$ws.Range("D11").PasteSpecial(-4122) | Out-Null
$ws.Range("D11").Value = "T_FL_FLOAT"
$ws.Range("K6").Copy() | Out-Null
$ws.Range("E11").PasteSpecial(-4122) | Out-Null
$ws.Range("E11").Value = "T_FL_FLOAT"
$ws.Range("K6").Copy() | Out-Null
$ws.Range("F11").PasteSpecial(-4122) | Out-Null
$ws.Range("F11").Value = "T_FL_FLOAT"
$ws.Range("K6").Copy() | Out-Null
$ws.Range("G11").PasteSpecial(-4122) | Out-Null
$ws.Range("G11").Value = "T_FL_FLOAT"
$ws.Range("K6").Copy() | Out-Null
$ws.Range("H11").PasteSpecial(-4122) | Out-Null
$ws.Range("H11").Value = "T_FL_FLOAT"
$ws.Range("K6").Copy() | Out-Null
$ws.Range("I11").PasteSpecial(-4122) | Out-Null
$ws.Range("I11").Value = "T_FL_FLOAT"

# row 12
$ws.Range("K1").Copy() | Out-Null
$ws.Range("A12").PasteSpecial(-4122) | Out-Null
$ws.Range("A12").Value = 10
$ws.Range("K5").Copy() | Out-Null
$ws.Range("B12").PasteSpecial(-4122) | Out-Null
$ws.Range("B12").Value = "T_FL_ERROR"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("C12").PasteSpecial(-4122) | Out-Null
$ws.Range("C12").Value = "T_FL_ERROR"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4122) | Out-Null
$ws.Range("D12").Value = "T_FL_ERROR"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("E12").PasteSpecial(-4122) | Out-Null
$ws.Range("E12").Value = "T_FL_ERROR"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("F12").PasteSpecial(-4122) | Out-Null
$ws.Range("F12").Value = "T_FL_ERROR"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("G12").PasteSpecial(-4122) | Out-Null
$ws.Range("G12").Value = "T_FL_ERROR"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("H12").PasteSpecial(-4122) | Out-Null
$ws.Range("H12").Value = "T_FL_ERROR"
$ws.Range("K3").Copy() | Out-Null
$ws.Range("I12").PasteSpecial(-4122) | Out-Null
$ws.Range("I12").Value = "T_FL_ERROR"

# --- Clear the scratch palette area.
$ws.Range("K1:K6").Clear() | Out-Null

# --- Sheet2 view: selection + dimension (dimension auto-updates from used range).
$ws.Activate()
$ws.Range("F21").Select() | Out-Null

# --- Sheet1 ("transitions"): zoom change 100 -> 85.
$ws1 = $wb.Worksheets.Item("transitions")
$ws1.Activate()
$excel.ActiveWindow.Zoom = 85
